$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# DATASETS sheet: remove stray empty H2 cell, add new row 4 for the
# consolidated xlsm dataset feeding the daily z-score indicator series.
# ---------------------------------------------------------------------------
$datasets = $wb.Worksheets.Item("DATASETS")

# H2 was an empty placeholder cell; clear it so it no longer exists.
$datasets.Range("H2").ClearContents()

$datasets.Range("A4").Value = "DAILY_CONSOLIDATED_XLSM"
$datasets.Range("B4").Value = "daily_series_wide"
$datasets.Range("C4").Value = "xlsm"
$datasets.Range("D4").Value = "E:\BacktestData\raw\consolidated.xlsm"
$datasets.Range("E4").Value = "append"
$datasets.Range("F4").Value = "Date"
$datasets.Range("G4").Value = "America/New_York"
$datasets.Range("H4").Value = "1D"
$datasets.Range("I4").Value = "eod_20_00_ny"
$datasets.Range("J4").Value = 0
$datasets.Range("K4").Value = "daily_series"
$datasets.Range("L4").Value = "year"
$datasets.Range("M4").Value = "Ingest columns: SpxCombined_pos,SpxSystematic_pos,SpxLS_pos,SpxMF_pos,SpxRetail_pos,Spx_NetOptionsPositioning,Spx_DlrGamma,EUshorts_pos,EUetf_pos,EUrp_pos,EUcta_pos,EULS_pos,EUMF_pos,EUComb_pos. Source file: consolidated.xlsm (append-only). Missing values allowed; coverage differs by series."

# ---------------------------------------------------------------------------
# INSTRUMENTS sheet: populate rows 2-15 with the new daily z-score series
# sourced from DAILY_CONSOLIDATED_XLSM.
# ---------------------------------------------------------------------------
$instruments = $wb.Worksheets.Item("INSTRUMENTS")

$spxNote = "Daily indicator from consolidated.xlsm; units=zscore; timing rule set in DATASETS.known_time_rule; per-series feature lags/transforms go in FEATURE_LIBRARY."

$rows = @(
    @{ Id = "SpxCombined_pos"; Currency = "USD"; Calendar = "NYSE" },
    @{ Id = "SpxSystematic_pos"; Currency = "USD"; Calendar = "NYSE" },
    @{ Id = "SpxLS_pos"; Currency = "USD"; Calendar = "NYSE" },
    @{ Id = "SpxMF_pos"; Currency = "USD"; Calendar = "NYSE" },
    @{ Id = "SpxRetail_pos"; Currency = "USD"; Calendar = "NYSE" },
    @{ Id = "Spx_NetOptionsPositioning"; Currency = "USD"; Calendar = "NYSE" },
    @{ Id = "Spx_DlrGamma"; Currency = "USD"; Calendar = "NYSE" },
    @{ Id = "EUshorts_pos"; Currency = "EUR"; Calendar = "EUREX" },
    @{ Id = "EUetf_pos"; Currency = "EUR"; Calendar = "EUREX" },
    @{ Id = "EUrp_pos"; Currency = "EUR"; Calendar = "EUREX" },
    @{ Id = "EUcta_pos"; Currency = "EUR"; Calendar = "EUREX" },
    @{ Id = "EULS_pos"; Currency = "EUR"; Calendar = "EUREX" },
    @{ Id = "EUMF_pos"; Currency = "EUR"; Calendar = "EUREX" },
    @{ Id = "EUComb_pos"; Currency = "EUR"; Calendar = "EUREX" }
)

$r = 2
foreach ($row in $rows) {
    $instruments.Cells.Item($r, 1).Value = $row.Id          # A instrument_id
    $instruments.Cells.Item($r, 2).Value = $row.Id          # B instrument_name
    $instruments.Cells.Item($r, 3).Value = "indicator"      # C instrument_type
    $instruments.Cells.Item($r, 4).Value = "DAILY_CONSOLIDATED_XLSM" # D prices_dataset_id
    $instruments.Cells.Item($r, 8).Value = $row.Id          # H close_col
    $instruments.Cells.Item($r, 12).Value = $row.Currency   # L currency
    $instruments.Cells.Item($r, 13).Value = $row.Calendar   # M calendar
    $instruments.Cells.Item($r, 16).Value = $spxNote        # P notes
    $r = $r + 1
}
